$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Rows("2:2").Insert()
$ws.Range("A2").Value = "T1"
